# "lock in current version" - rewrite the cue sequence table (word / image / category)
# with the new locked-in word list and matching face/flower image assignments.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("wenden", "none", "none"),
    @("leiten", "flower/flower031.jpg", "flower"),
    @("segnen", "face/face031.jpg", "face"),
    @("rufen", "none", "none"),
    @("sehen", "flower/flower027.jpg", "flower"),
    @("wandern", "face/face012.jpg", "face"),
    @("enden", "none", "none"),
    @("nerven", "flower/flower010.jpg", "flower"),
    @("erben", "flower/flower012.jpg", "flower"),
    @("opfern", "none", "none"),
    @("faulen", "flower/flower024.jpg", "flower"),
    @("jubeln", "face/face009.jpg", "face"),
    @("schalten", "none", "none"),
    @("mauern", "face/face015.jpg", "face"),
    @("spielen", "flower/flower017.jpg", "flower"),
    @("klagen", "none", "none"),
    @("ändern", "flower/flower006.jpg", "flower"),
    @("wagen", "flower/flower015.jpg", "flower"),
    @("sparen", "none", "none"),
    @("bieten", "face/face003.jpg", "face"),
    @("lassen", "face/face006.jpg", "face"),
    @("hören", "none", "none"),
    @("heben", "face/face016.jpg", "face"),
    @("leisten", "flower/flower033.jpg", "flower"),
    @("drohen", "none", "none"),
    @("schütteln", "flower/flower018.jpg", "flower"),
    @("binden", "flower/flower021.jpg", "flower"),
    @("orten", "none", "none"),
    @("antun", "face/face023.jpg", "face"),
    @("kümmern", "face/face020.jpg", "face"),
    @("weigern", "none", "none"),
    @("schweben", "flower/flower020.jpg", "flower"),
    @("gleichen", "flower/flower028.jpg", "flower"),
    @("ächzen", "none", "none"),
    @("bitten", "face/face005.jpg", "face"),
    @("schaden", "face/face018.jpg", "face"),
    @("kosten", "none", "none"),
    @("fügen", "flower/flower022.jpg", "flower"),
    @("ruhen", "face/face011.jpg", "face"),
    @("dauern", "none", "none"),
    @("stärken", "face/face010.jpg", "face"),
    @("danken", "face/face004.jpg", "face"),
    @("stören", "none", "none"),
    @("münzen", "flower/flower002.jpg", "flower"),
    @("platzen", "face/face026.jpg", "face"),
    @("bremsen", "none", "none"),
    @("klingen", "face/face017.jpg", "face"),
    @("lächeln", "flower/flower001.jpg", "flower")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}